$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 through 14 (old rows that get merged/removed)
$ws.Rows("6:14").Delete()

# Update A2:A5 with the new consolidated values
$ws.Range("A2").Value = "('Elspeth, Knight-Errant Emblem', ['Emblem — Elspeth', 'Artifacts, creatures, enchantments, and lands you control have indestructible.'])"
$ws.Range("A3").Value = "('Myr', ['Token Artifact Creature — Myr', '1/1'])"
$ws.Range("A4").Value = "('Soldier', ['Token Creature — Soldier', '1/1'])"
$ws.Range("A5").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
